# Arshdeep Singh.xlsx - add matchNo column + a new scraped match row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab to match the player name.
$ws.Name = "Arshdeep Singh"

# Insert a brand new column A ("matchNo"); this shifts the existing
# teamName..result columns (previously A:L) right by one, to B:M.
$ws.Columns("A:A").Insert()

# The sheet stores every value as text (see the numberStoredAsText
# ignoredError), so force the whole used range to text formatting before
# writing any numeric-looking strings (e.g. "1", "0", "100.00").
$ws.Range("A1:M3").NumberFormat = "@"

# --- Header row ---
$ws.Range("A1").Value = "matchNo"

# --- Row 2: fill in the new matchNo value for the existing match, and
#     correct the balls-faced figure that was re-scraped (0 -> 1). ---
$ws.Range("A2").Value = "21st"
$ws.Range("F2").Value = "1"

# --- Row 3: brand new match scraped for Arshdeep Singh. ---
$ws.Range("A3").Value = "14th"
$ws.Range("B3").Value = "Punjab Kings"
$ws.Range("C3").Value = "Arshdeep Singh"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "1"
$ws.Range("F3").Value = "2"
$ws.Range("G3").Value = "0"
$ws.Range("H3").Value = "0"
$ws.Range("I3").Value = "50.00"
$ws.Range("J3").Value = "Sunrisers Hyderabad"
$ws.Range("K3").Value = "Chennai"
$ws.Range("L3").Value = "April 21"
$ws.Range("M3").Value = "Sunrisers won by 9 wickets (with 8 balls remaining)"
